# Updates the cryptos price list (Coin / Link / Price / Volume(1h)) to the
# latest scraped values. Price values in column D are forced to text
# (leading apostrophe) so strings like "23.910.47" or "0.3890" are not
# reinterpreted by Excel as numbers/dates, then the cell style is restored
# to "Normal" so no stray number-format is left visible on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.910.47"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "'1.653.26"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.3890"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "'51.67"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").Value = "'1.354"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("D11").Value = "'1.002"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "'0.08487"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "'24.11"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").Value = "'7.074"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").Value = "'8.104"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "'0.00001320"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "'1.651.79"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "'94.32"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'0.07004"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'19.72"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("D21").Value = "'6.982"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'23.905.90"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("D27").Value = "'22.13"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").Value = "'154.11"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").Value = "'5.419"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Value = "'138.27"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").Value = "'7.898"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'2.509"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "'1.838.08"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").Value = "'1.027"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("D35").Value = "'0.08177"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "'6.713"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("D37").Value = "'0.02921"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("D38").Value = "'10.87"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("D39").Value = "'0.2684"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").Value = "'0.09149"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'13.78"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.7589"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "'1.428"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("D44").Value = "'16.56"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "'0.6969"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").Value = "'2.466"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").Value = "'4.103"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'1.0000"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "'0.08308"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "'134.48"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "'1.234"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = "  -2.90%  "